$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '63.795.01'
$ws.Range("E2").Value = '  -0.70%  '

# Row 3
$ws.Range("D3").Value = '3.436.07'
$ws.Range("E3").Value = '  -1.60%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.31'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.03%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.09'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.02%  '

# Row 7
$ws.Range("E7").Value = '  +0.00%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.479'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.64%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.62'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.22%  '

# Row 10
$ws.Range("E10").Value = '  -1.48%  '

# Row 11
$ws.Range("E11").Value = '  -0.95%  '

# Row 12
$ws.Range("D12").Value = '4.035.81'
$ws.Range("E12").Value = '  -1.22%  '

# Row 13
$ws.Range("E13").Value = '  -0.41%  '

# Row 14
$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000176'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.01%  '

# Row 15
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '3.442.48'
$ws.Range("E15").Value = '  -1.82%  '

# Row 16
$ws.Range("D16").Value = '63.694.94'
$ws.Range("E16").Value = '  -0.99%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '24.82'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.75%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.64'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.86%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.29'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.26%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '383.76'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.62%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.561'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.65%  '

# Row 23
$ws.Range("D23").Value = '3.580.14'
$ws.Range("E23").Value = '  -1.42%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.97%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.996'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.39%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.37'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.34%  '

# Row 27
$ws.Range("E27").Value = '  -4.66%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.996'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.19%  '

# Row 29
$ws.Range("E29").Value = '  -1.01%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.02'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.27%  '

# Row 31
$ws.Range("E31").Value = '  +2.55%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.89'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.96%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.42'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.13%  '

# Row 34
$ws.Range("D34").Value = '3.466.51'
$ws.Range("E34").Value = '  -1.35%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '22.80'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.54%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.16'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.64%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.71'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.34%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '163.57'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.84%  '

# Row 40
$ws.Range("E40").Value = '  -4.27%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0769'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.18%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.791'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.75%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.10%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.25'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.35%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.32'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.49%  '

# Row 46
$ws.Range("E46").Value = '  -2.84%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.47'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.14%  '

# Row 48
$ws.Range("E48").Value = '  -4.79%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.68'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.06%  '

# Row 50
$ws.Range("D50").Value = '2.312.72'
$ws.Range("E50").Value = '  -6.00%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.885'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.87%  '
